$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-14 (B, C, D, E columns) with corrected schedule data
$ws.Range("B3").Value = "03:30"
$ws.Range("C3").Value = "ΕΛΕΝΗ"
$ws.Range("D3").Value = "06:00"
$ws.Range("E3").Value = "ΕΛΕΝΗ"

$ws.Range("B4").Value = "05:30"
$ws.Range("C4").Value = "ΙΩΝΑΣ"
$ws.Range("D4").Value = "07:30"
$ws.Range("E4").Value = "ΑΓΙΑ ΕΙΡΗΝΗ"

$ws.Range("B5").Value = "09:30"
$ws.Range("C5").Value = "ΕΛΕΝΗ"
$ws.Range("D5").Value = "08:45"
$ws.Range("E5").Value = "ΕΡΜΗΣ"

$ws.Range("B6").Value = "10:00"
$ws.Range("C6").Value = "ΝΑΝΤΗ"
$ws.Range("D6").Value = "10:30"
$ws.Range("E6").Value = "ΙΩΝΑΣ"

$ws.Range("B7").Value = "10:45"
$ws.Range("C7").Value = "ΕΡΜΗΣ"
$ws.Range("D7").Value = "12:15"
$ws.Range("E7").Value = "ΕΛΕΝΗ"

$ws.Range("B8").Value = "12:00"
$ws.Range("C8").Value = "ΑΓΙΑ ΕΙΡΗΝΗ"
$ws.Range("D8").Value = "13:00"
$ws.Range("E8").Value = "ΑΓΙΟΣ ΣΠΥΡΙΔΩΝ"

$ws.Range("B9").Value = "12:30"
$ws.Range("C9").Value = "ΙΩΝΑΣ"
$ws.Range("D9").Value = "13:45"
$ws.Range("E9").Value = "ΕΡΜΗΣ"

$ws.Range("B10").Value = "15:45"
$ws.Range("C10").Value = "ΕΡΜΗΣ"
$ws.Range("D10").Value = "14:45"
$ws.Range("E10").Value = "ΙΩΝΑΣ"

$ws.Range("B11").Value = "18:00"
$ws.Range("C11").Value = "ΑΓΙΑ ΕΙΡΗΝΗ"
$ws.Range("D11").Value = "15:45"
$ws.Range("E11").Value = "ΑΓΙΑ ΕΙΡΗΝΗ"

$ws.Range("B12").Value = "18:30"
$ws.Range("C12").Value = "ΑΓΙΟΣ ΣΠΥΡΙΔΩΝ"
$ws.Range("D12").Value = "17:00"
$ws.Range("E12").Value = "ΝΑΝΤΗ"

$ws.Range("B13").Value = "19:15"
$ws.Range("C13").Value = "ΝΑΝΤΗ"
$ws.Range("D13").Value = "17:45"
$ws.Range("E13").Value = "ΕΡΜΗΣ"

$ws.Range("B14").Value = "19:45"
$ws.Range("C14").Value = "ΕΡΜΗΣ"
$ws.Range("D14").Value = "20:45"
$ws.Range("E14").Value = "ΑΓΙΟΣ ΣΠΥΡΙΔΩΝ"

# Row 15 is no longer part of the schedule; remove it entirely so the
# dimension shrinks from B2:E15 to B2:E14
$ws.Rows("15:15").Delete()
